# Auto-generated edit script: update Leve profit-tracking values per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 86
$ws.Range("H86").Value = 3017.0833
$ws.Range("I86").Value = 3182.2727
$ws.Range("J86").Value = 1200
$ws.Range("K86").Value = 3182.2727
$ws.Range("L86").Value = 1200
$ws.Range("M86").Value = -2059.2727
$ws.Range("N86").Value = -3446
# Row 89
$ws.Range("H89").Value = 3017.0833
$ws.Range("I89").Value = 3182.2727
$ws.Range("J89").Value = 1200
$ws.Range("K89").Value = 15911.3635
$ws.Range("L89").Value = 6000
$ws.Range("M89").Value = -10295.3635
$ws.Range("N89").Value = -17232

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 27114.906
$ws.Range("I32").Value = 29115.227
$ws.Range("J32").Value = 17477
$ws.Range("K32").Value = 29115.227
$ws.Range("L32").Value = 17477
$ws.Range("M32").Value = -28828.227
$ws.Range("N32").Value = -18051
# Row 64
$ws.Range("H64").Value = 29750
$ws.Range("J64").Value = 29750
$ws.Range("L64").Value = 29750
$ws.Range("N64").Value = -30246
# Row 67
$ws.Range("H67").Value = 29750
$ws.Range("J67").Value = 29750
$ws.Range("L67").Value = 29750
$ws.Range("N67").Value = -31466
# Row 74
$ws.Range("H74").Value = 1756.7858
$ws.Range("I74").Value = 990.8570999999999
$ws.Range("J74").Value = 2522.7144
$ws.Range("K74").Value = 990.8570999999999
$ws.Range("L74").Value = 2522.7144
$ws.Range("M74").Value = -116.8570999999999
$ws.Range("N74").Value = -4270.7144
# Row 77
$ws.Range("H77").Value = 1756.7858
$ws.Range("I77").Value = 990.8570999999999
$ws.Range("J77").Value = 2522.7144
$ws.Range("K77").Value = 4954.2855
$ws.Range("L77").Value = 12613.572
$ws.Range("M77").Value = -586.2855
$ws.Range("N77").Value = -21349.572
# Row 118
$ws.Range("H118").Value = 49537
$ws.Range("J118").Value = 49537
$ws.Range("L118").Value = 49537
$ws.Range("N118").Value = -52851
# Row 120
$ws.Range("H120").Value = 42197.332
$ws.Range("J120").Value = 42197.332
$ws.Range("L120").Value = 42197.332
$ws.Range("N120").Value = -51873.332
# Row 123
$ws.Range("H123").Value = 43108.5
$ws.Range("J123").Value = 43108.5
$ws.Range("L123").Value = 43108.5
$ws.Range("N123").Value = -52908.5
# Row 132
$ws.Range("H132").Value = 35716436
$ws.Range("I132").Value = 50001400
$ws.Range("J132").Value = 4024
$ws.Range("K132").Value = 150004200
$ws.Range("L132").Value = 12072
$ws.Range("M132").Value = -150001670
$ws.Range("N132").Value = -17132
# Row 134
$ws.Range("H134").Value = 37228.062
$ws.Range("J134").Value = 37228.062
$ws.Range("L134").Value = 37228.062
$ws.Range("N134").Value = -47368.062
# Row 137
$ws.Range("H137").Value = 23065.75
$ws.Range("J137").Value = 25098.182
$ws.Range("L137").Value = 25098.182
$ws.Range("N137").Value = -35298.182
# Row 138
$ws.Range("H138").Value = 47518.332
$ws.Range("J138").Value = 47518.332
$ws.Range("L138").Value = 47518.332
$ws.Range("N138").Value = -57798.332

$ws = $wb.Worksheets.Item("BSM")
# Row 117
$ws.Range("H117").Value = 44999
$ws.Range("J117").Value = 44999
$ws.Range("L117").Value = 44999
$ws.Range("N117").Value = -54177
# Row 119
$ws.Range("H119").Value = 47992
$ws.Range("J119").Value = 47992
$ws.Range("L119").Value = 47992
$ws.Range("N119").Value = -57668
# Row 120
$ws.Range("H120").Value = 48753
$ws.Range("J120").Value = 48753
$ws.Range("L120").Value = 48753
$ws.Range("N120").Value = -58429
# Row 122
$ws.Range("H122").Value = 40725.332
$ws.Range("J122").Value = 40725.332
$ws.Range("L122").Value = 40725.332
$ws.Range("N122").Value = -50525.332
# Row 132
$ws.Range("H132").Value = 38796.875
$ws.Range("J132").Value = 38796.875
$ws.Range("L132").Value = 38796.875
$ws.Range("N132").Value = -48916.875
# Row 139
$ws.Range("H139").Value = 60833
$ws.Range("J139").Value = 60833
$ws.Range("L139").Value = 60833
$ws.Range("N139").Value = -71113

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 217651.9
$ws.Range("J31").Value = 243068.02
$ws.Range("L31").Value = 243068.02
$ws.Range("N31").Value = -243658.02
# Row 34
$ws.Range("H34").Value = 217651.9
$ws.Range("J34").Value = 243068.02
$ws.Range("L34").Value = 243068.02
$ws.Range("N34").Value = -243472.02
# Row 37
$ws.Range("H37").Value = 41999
$ws.Range("J37").Value = 41999
$ws.Range("L37").Value = 41999
$ws.Range("N37").Value = -42213
# Row 58
$ws.Range("H58").Value = 1925.4688
$ws.Range("I58").Value = 1691.25
$ws.Range("J58").Value = 2628.125
$ws.Range("K58").Value = 1691.25
$ws.Range("L58").Value = 2628.125
$ws.Range("M58").Value = -1488.25
$ws.Range("N58").Value = -3034.125
# Row 118
$ws.Range("H118").Value = 44534
$ws.Range("J118").Value = 44534
$ws.Range("L118").Value = 44534
$ws.Range("N118").Value = -47848
# Row 121
$ws.Range("H121").Value = 32982.25
$ws.Range("J121").Value = 32982.25
$ws.Range("L121").Value = 32982.25
$ws.Range("N121").Value = -35602.25
# Row 133
$ws.Range("H133").Value = 24983.166
$ws.Range("J133").Value = 24983.166
$ws.Range("L133").Value = 24983.166
$ws.Range("N133").Value = -30043.166
# Row 136
$ws.Range("H136").Value = 1925.4688
$ws.Range("I136").Value = 1691.25
$ws.Range("J136").Value = 2628.125
$ws.Range("K136").Value = 5073.75
$ws.Range("L136").Value = 7884.375
$ws.Range("M136").Value = -2523.75
$ws.Range("N136").Value = -12984.375

$ws = $wb.Worksheets.Item("CUL")
# Row 25
$ws.Range("H25").Value = 1000
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
# Row 30
$ws.Range("H30").Value = 1000
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
# Row 97
$ws.Range("H97").Value = 1000
$ws.Range("J97").Value = 1000
$ws.Range("L97").Value = 3000
$ws.Range("N97").Value = -3992
# Row 123
$ws.Range("H123").Value = 7433.3335
$ws.Range("J123").Value = 2920
$ws.Range("L123").Value = 8760
$ws.Range("N123").Value = -13660

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 29414498
$ws.Range("I7").Value = 45456810
$ws.Range("J7").Value = 3600.4167
$ws.Range("K7").Value = 45456810
$ws.Range("L7").Value = 3600.4167
$ws.Range("M7").Value = -45456698
$ws.Range("N7").Value = -3824.4167
# Row 12
$ws.Range("H12").Value = 35224
$ws.Range("I12").Value = 444
$ws.Range("J12").Value = 70004
$ws.Range("K12").Value = 444
$ws.Range("L12").Value = 70004
$ws.Range("M12").Value = -274
$ws.Range("N12").Value = -70344
# Row 97
$ws.Range("H97").Value = 34996
$ws.Range("J97").Value = 34996
$ws.Range("L97").Value = 34996
$ws.Range("N97").Value = -36978
# Row 104
$ws.Range("H104").Value = 17280
$ws.Range("J104").Value = 17280
$ws.Range("L104").Value = 17280
$ws.Range("N104").Value = -24268
# Row 121
$ws.Range("H121").Value = 40944
$ws.Range("J121").Value = 40944
$ws.Range("L121").Value = 40944
$ws.Range("N121").Value = -44438
# Row 126
$ws.Range("H126").Value = 29414498
$ws.Range("I126").Value = 45456810
$ws.Range("J126").Value = 3600.4167
$ws.Range("K126").Value = 136370430
$ws.Range("L126").Value = 10801.2501
$ws.Range("M126").Value = -136367960
$ws.Range("N126").Value = -15741.2501
# Row 132
$ws.Range("H132").Value = 3087.0715
$ws.Range("I132").Value = 1383.35
$ws.Range("J132").Value = 4635.909
$ws.Range("K132").Value = 4150.049999999999
$ws.Range("L132").Value = 13907.727
$ws.Range("M132").Value = -1620.049999999999
$ws.Range("N132").Value = -18967.727
# Row 136
$ws.Range("H136").Value = 1594.5625
$ws.Range("I136").Value = 1066.3334
$ws.Range("J136").Value = 3179.25
$ws.Range("K136").Value = 3199.0002
$ws.Range("L136").Value = 9537.75
$ws.Range("M136").Value = -649.0001999999999
$ws.Range("N136").Value = -14637.75
# Row 137
$ws.Range("H137").Value = 27234
$ws.Range("J137").Value = 27234
$ws.Range("L137").Value = 27234
$ws.Range("N137").Value = -37434
# Row 139
$ws.Range("H139").Value = 36634.617
$ws.Range("J139").Value = 36634.617
$ws.Range("L139").Value = 36634.617
$ws.Range("N139").Value = -46914.617

$ws = $wb.Worksheets.Item("WVR")
# Row 15
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
# Row 16
$ws.Range("H16").Value = 41779.6
$ws.Range("J16").Value = 41779.6
$ws.Range("L16").Value = 41779.6
$ws.Range("N16").Value = -42363.6
# Row 18
$ws.Range("H18").Value = 9828.333000000001
$ws.Range("I18").Value = 8970
$ws.Range("J18").Value = 10000
$ws.Range("K18").Value = 8970
$ws.Range("L18").Value = 10000
$ws.Range("M18").Value = -8797
$ws.Range("N18").Value = -10346
# Row 120
$ws.Range("H120").Value = 45412
$ws.Range("J120").Value = 45412
$ws.Range("L120").Value = 45412
$ws.Range("N120").Value = -55088
# Row 136
$ws.Range("H136").Value = 16953.438
$ws.Range("I136").Value = 35166.035
$ws.Range("K136").Value = 105498.105
$ws.Range("M136").Value = -102948.105
# Row 139
$ws.Range("H139").Value = 20676.875
$ws.Range("J139").Value = 20676.875
$ws.Range("L139").Value = 20676.875
$ws.Range("N139").Value = -30956.875
